$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.974.04"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.427.49"
$ws.Range("E3").Value = "  +5.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.63"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.16"
$ws.Range("E6").Value = "  +10.15%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.424.85"
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +6.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.23"
$ws.Range("E14").Value = "  +13.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.861.79"
$ws.Range("E15").Value = "  +5.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.857.19"
$ws.Range("E16").Value = "  +5.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +8.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.427.04"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.21"
$ws.Range("E19").Value = "  +7.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.43"
$ws.Range("E20").Value = "  +9.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +5.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.17"
$ws.Range("E24").Value = "  +3.45%  "
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.52"
$ws.Range("E27").Value = "  +13.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  +5.67%  "
$ws.Range("E29").Value = "  +11.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0794"
$ws.Range("E30").Value = "  +11.12%  "
$ws.Range("E31").Value = "  +7.14%  "
$ws.Range("E32").Value = "  +14.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "174.33"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +11.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.398"
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.69"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "369.74"
$ws.Range("E37").Value = "  +18.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.48"
$ws.Range("E38").Value = "  +12.40%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +13.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.00"
$ws.Range("E42").Value = "  +6.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.77"
$ws.Range("E43").Value = "  +9.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.69"
$ws.Range("E44").Value = "  +8.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.58"
$ws.Range("E45").Value = "  +11.56%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0957"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.591"
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +6.56%  "
$ws.Range("E49").Value = "  +6.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.79"
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("E51").Value = "  +16.49%  "
